$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Fill in the "完成情况" (completion status) column for the previous
#    week's block (rows 153-158) with "已完成".
# ---------------------------------------------------------------------------
$ws.Range("C153").Value = "已完成"
$ws.Range("C154").Value = "已完成"
$ws.Range("C155").Value = "已完成"
$ws.Range("C156").Value = "已完成"
$ws.Range("C157").Value = "已完成"
$ws.Range("C158").Value = "已完成"

# ---------------------------------------------------------------------------
# 2. Fill in the summary text for that same week (was just the "总结：" label).
# ---------------------------------------------------------------------------
$ws.Range("A159").Value = "总结：本此对app进一步功能完善，对websocket聊天功能代码进行了修改和优化。聊天功能上计划新增收发图片功能，下一次计划正式实施。"

# ---------------------------------------------------------------------------
# 3. Append a brand-new weekly block (rows 161-170), mirroring the layout of
#    the previous block (rows 151-160: date row, header row, six member
#    rows, then a two-row merged summary placeholder).
#    Merge the new header/summary ranges *before* copying the formatting
#    over, so the merge operation doesn't strip/rewrite the interior
#    borders of a range that already carries data.
# ---------------------------------------------------------------------------
$ws.Range("A161:D161").Merge()
$ws.Range("A169:D170").Merge()

$ws.Range("A151:D160").Copy()
$ws.Range("A161").PasteSpecial(-4122)

# Date heading with the weekday portion bolded (matches the existing
# pattern used for every other weekly date heading in the sheet).
$ws.Range("A161").Value = "日期：2018.11.12 第十一周周一"
$ws.Range("A161").Characters(14, 7).Font.Bold = $true

# Column headers (组员 / 计划内容 / 完成情况 / 备注).
$ws.Range("A162").Value = "组员"
$ws.Range("B162").Value = "计划内容"
$ws.Range("C162").Value = "完成情况"
$ws.Range("D162").Value = "备注"

# Member rows.
$ws.Range("A163").Value = "李光洪"
$ws.Range("B163").Value = "web端管理员端群组管理界面编码"

$ws.Range("A164").Value = "吴彤林"
$ws.Range("B164").Value = "app端创建好友分组功能实现与接口对接"

$ws.Range("A165").Value = "劳汉文"
$ws.Range("B165").Value = "app端移动好友分组功能实现与接口对接"

$ws.Range("A166").Value = "方嘉耀"
$ws.Range("B166").Borders.LineStyle = -4142
$ws.Range("B166").Value = "app端删除好友分组功能实现与接口对接"

$ws.Range("A167").Value = "成世靖"
$ws.Range("B167").Value = "app端新增好友功能、删除好友功能实现与接口对接"

$ws.Range("A168").Value = "丰浩"
$ws.Range("B168").Value = "web端群组界面协助编码"

# The "完成情况" cell for row 164 was left fully blank (no cell format
# either) in the source edit, unlike its siblings.
$ws.Range("C164").Clear()

# Summary placeholder row for the new block ("总结：", content pending).
$ws.Range("A169").Value = "总结："

# ---------------------------------------------------------------------------
# 4. Leave the selection where the editor ended up after typing the new
#    block, scrolled down so the new rows are visible.
# ---------------------------------------------------------------------------
$ws.Range("B168").Select()
